# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to several Leve-profit rows across 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1786.1063
$ws.Range("I15").Value = 1786.1063
$ws.Range("K15").Value = 5358.3189
$ws.Range("M15").Value = -5189.3189
$ws.Range("H41").Value = 1410.5333
$ws.Range("I41").Value = 820.4286
$ws.Range("K41").Value = 820.4286
$ws.Range("M41").Value = -380.4286
$ws.Range("H43").Value = 1713.8572
$ws.Range("I43").Value = 1866
$ws.Range("K43").Value = 1866
$ws.Range("M43").Value = -1797
$ws.Range("H49").Value = 4438.3335
$ws.Range("I49").Value = 1999
$ws.Range("J49").Value = 5658
$ws.Range("K49").Value = 5997
$ws.Range("L49").Value = 16974
$ws.Range("M49").Value = -5861
$ws.Range("N49").Value = -17246
$ws.Range("H53").Value = 672.7083
$ws.Range("I53").Value = 363.26666
$ws.Range("K53").Value = 363.26666
$ws.Range("M53").Value = 273.73334
$ws.Range("H62").Value = 7201.0586
$ws.Range("I62").Value = 3141.8
$ws.Range("J62").Value = 13000
$ws.Range("K62").Value = 3141.8
$ws.Range("L62").Value = 13000
$ws.Range("M62").Value = -2517.8
$ws.Range("N62").Value = -14248
$ws.Range("H65").Value = 7201.0586
$ws.Range("I65").Value = 3141.8
$ws.Range("J65").Value = 13000
$ws.Range("K65").Value = 15709
$ws.Range("L65").Value = 65000
$ws.Range("M65").Value = -12589
$ws.Range("N65").Value = -71240
$ws.Range("H70").Value = 8332.200000000001
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 8332.200000000001
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 24996.6
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -25536.6
$ws.Range("H73").Value = 8332.200000000001
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 8332.200000000001
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 24996.6
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -26868.6
$ws.Range("H132").Value = 12281.277
$ws.Range("I132").Value = 12281.277
$ws.Range("K132").Value = 36843.831
$ws.Range("M132").Value = -34313.831

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7799.933
$ws.Range("I32").Value = 6714.2144
$ws.Range("K32").Value = 6714.2144
$ws.Range("M32").Value = -6427.2144
$ws.Range("H45").Value = 3063.625
$ws.Range("I45").Value = 2158.75
$ws.Range("K45").Value = 2158.75
$ws.Range("M45").Value = -1781.75

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 600
$ws.Range("I10").Value = 600
$ws.Range("K10").Value = 600
$ws.Range("M10").Value = -460
$ws.Range("H94").Value = 430.3846
$ws.Range("I94").Value = 332.8889
$ws.Range("J94").Value = 649.75
$ws.Range("K94").Value = 332.8889
$ws.Range("L94").Value = 649.75
$ws.Range("M94").Value = 118.1111
$ws.Range("N94").Value = -1551.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5217.2
$ws.Range("I58").Value = 2500
$ws.Range("K58").Value = 2500
$ws.Range("M58").Value = -2297
$ws.Range("H134").Value = 2065.6
$ws.Range("I134").Value = 2106.2222
$ws.Range("K134").Value = 6318.6666
$ws.Range("M134").Value = -3783.6666
$ws.Range("H136").Value = 5217.2
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 143142900
$ws.Range("J4").Value = 81.333336
$ws.Range("L4").Value = 244.000008
$ws.Range("N4").Value = -468.000008
$ws.Range("H44").Value = 716.3570999999999
$ws.Range("I44").Value = 170
$ws.Range("K44").Value = 510
$ws.Range("M44").Value = -112

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 37575
$ws.Range("J46").Value = 37575
$ws.Range("L46").Value = 37575
$ws.Range("N46").Value = -37887
$ws.Range("H80").Value = 2228.2856
$ws.Range("I80").Value = 2283
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 2283
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -1285
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 2228.2856
$ws.Range("I83").Value = 2283
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 11415
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -6423
$ws.Range("N83").Value = -19484
$ws.Range("H132").Value = 105110.1
$ws.Range("I132").Value = 253426.25
$ws.Range("K132").Value = 760278.75
$ws.Range("M132").Value = -757748.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 911.3333
$ws.Range("I16").Value = 920.25
$ws.Range("J16").Value = 893.5
$ws.Range("K16").Value = 920.25
$ws.Range("L16").Value = 893.5
$ws.Range("M16").Value = -750.25
$ws.Range("N16").Value = -1233.5
$ws.Range("H40").Value = 4142.75
$ws.Range("I40").Value = 3398
$ws.Range("J40").Value = 4887.5
$ws.Range("K40").Value = 3398
$ws.Range("L40").Value = 4887.5
$ws.Range("M40").Value = -3262
$ws.Range("N40").Value = -5159.5
$ws.Range("H68").Value = 8249.875
$ws.Range("I68").Value = 6666.6665
$ws.Range("J68").Value = 9199.799999999999
$ws.Range("K68").Value = 6666.6665
$ws.Range("L68").Value = 9199.799999999999
$ws.Range("M68").Value = -5917.6665
$ws.Range("N68").Value = -10697.8
$ws.Range("H71").Value = 8249.875
$ws.Range("I71").Value = 6666.6665
$ws.Range("J71").Value = 9199.799999999999
$ws.Range("K71").Value = 33333.3325
$ws.Range("L71").Value = 45999
$ws.Range("M71").Value = -29589.3325
$ws.Range("N71").Value = -53487
$ws.Range("H74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("K74").Value = 20000
$ws.Range("M74").Value = -19002
$ws.Range("H77").Value = 20000
$ws.Range("I77").Value = 20000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55008

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9483.166999999999
$ws.Range("I62").Value = 8359.799999999999
$ws.Range("J62").Value = 10285.571
$ws.Range("K62").Value = 8359.799999999999
$ws.Range("L62").Value = 10285.571
$ws.Range("M62").Value = -7735.799999999999
$ws.Range("N62").Value = -11533.571
$ws.Range("H65").Value = 9483.166999999999
$ws.Range("I65").Value = 8359.799999999999
$ws.Range("J65").Value = 10285.571
$ws.Range("K65").Value = 41799
$ws.Range("L65").Value = 51427.855
$ws.Range("M65").Value = -38679
$ws.Range("N65").Value = -57667.855
$ws.Range("H81").Value = 2199.6
$ws.Range("I81").Value = 999
$ws.Range("K81").Value = 1998
$ws.Range("M81").Value = -937
$ws.Range("H84").Value = 2199.6
$ws.Range("I84").Value = 999
$ws.Range("K84").Value = 9990
$ws.Range("M84").Value = -4686
$ws.Range("H98").Value = 35590
$ws.Range("J98").Value = 35590
$ws.Range("L98").Value = 35590
$ws.Range("N98").Value = -41580
